$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Unit" mapping column (B) gains a new "LCU, % of GDP, SA" unit:
# row 10 changes from the generic "% of GDP" to "LCU, % of GDP", a new
# row 11 "LCU, % of GDP, SA" is inserted, and the remaining "*, SA" units
# that used to occupy rows 11-16 shift down one row to rows 12-17.
$ws.Range("B10").Value = "LCU, % of GDP"
$ws.Range("B11").Value = "LCU, % of GDP, SA"
$ws.Range("B12").Value = "LCU, SA"
$ws.Range("B13").Value = "USD, SA"
$ws.Range("B14").Value = "LCU, % YoY, SA"
$ws.Range("B15").Value = "USD, % YoY, SA"
$ws.Range("B16").Value = "LCU, % QoQ, SA"
$ws.Range("B17").Value = "USD, % QoQ, SA"

# Move the active selection to where the author left off working (H19).
$ws.Range("H19").Select() | Out-Null
